$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.78947622181619
$ws.Range("C2").Value = 11.91735769457311
$ws.Range("E2").Value = 16.64639026427051
$ws.Range("F2").Value = 16.86991607391245
$ws.Range("G2").Value = 13.44684530066363
$ws.Range("H2").Value = 10.46931263882562
$ws.Range("O2").Value = 13.8656360093813

$ws.Range("B3").Value = 13.91390094346919
$ws.Range("C3").Value = 11.35696793297774
$ws.Range("E3").Value = 15.69293043590993
$ws.Range("F3").Value = 15.89584955866815
$ws.Range("G3").Value = 13.62221876042845
$ws.Range("H3").Value = 10.55183101261696
$ws.Range("O3").Value = 14.0152184301655

$ws.Range("B4").Value = 13.34598600319462
$ws.Range("C4").Value = 10.99756377274626
$ws.Range("E4").Value = 15.08172198445245
$ws.Range("F4").Value = 15.26997757108489
$ws.Range("G4").Value = 13.74444445141016
$ws.Range("H4").Value = 10.60545738609867
$ws.Range("O4").Value = 14.1133510405744

$ws.Range("B5").Value = 13.10703603479725
$ws.Range("C5").Value = 10.84739897228758
$ws.Range("E5").Value = 14.82643165056157
$ws.Range("F5").Value = 15.008197319934
$ws.Range("G5").Value = 13.79782721113864
$ws.Range("H5").Value = 10.62805429764795
$ws.Range("O5").Value = 14.15491284371824

$ws.Range("B6").Value = 13.06690785370488
$ws.Range("C6").Value = 10.82224528039201
$ws.Range("E6").Value = 14.78367412364066
$ws.Range("F6").Value = 14.96433081551589
$ws.Range("G6").Value = 13.80690486966837
$ws.Range("H6").Value = 10.63185140305509
$ws.Range("O6").Value = 14.16190882778707

$ws.Range("B7").Value = 13.34279373615794
$ws.Range("C7").Value = 10.99555338411016
$ws.Range("E7").Value = 15.07830383703977
$ws.Range("F7").Value = 15.26647399323133
$ws.Range("G7").Value = 13.74515002517509
$ws.Range("H7").Value = 10.60575912529567
$ws.Range("O7").Value = 14.11390520483222

$ws.Range("B8").Value = 14.4939214749363
$ws.Range("C8").Value = 11.72739721384023
$ws.Range("E8").Value = 16.32313114996975
$ws.Range("F8").Value = 16.5399640634477
$ws.Range("G8").Value = 13.50424528789133
$ws.Range("H8").Value = 10.49715034438591
$ws.Range("O8").Value = 13.9159011693209

$ws.Range("B9").Value = 16.50750969901804
$ws.Range("C9").Value = 13.03558586613095
$ws.Range("E9").Value = 18.68393690344874
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 13.15093055539266
$ws.Range("H9").Value = 10.30767810987004
$ws.Range("O9").Value = 13.5779442121903

$ws.Range("B10").Value = 17.83516559957563
$ws.Range("C10").Value = 13.91348294539697
$ws.Range("E10").Value = 20.34426609302905
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 12.96891464568268
$ws.Range("H10").Value = 10.18284132089146
$ws.Range("O10").Value = 13.3609483576492

$ws.Range("B11").Value = 18.40582767721765
$ws.Range("C11").Value = 14.29388445441589
$ws.Range("E11").Value = 21.0569017443305
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 12.90397486938251
$ws.Range("H11").Value = 10.1291796216444
$ws.Range("O11").Value = 13.26917121561851

$ws.Range("B12").Value = 18.61711567661016
$ws.Range("C12").Value = 14.43514826719364
$ws.Range("E12").Value = 21.32066146033916
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 12.88203176778028
$ws.Range("H12").Value = 10.10931013226337
$ws.Range("O12").Value = 13.23542766316013

$ws.Range("B13").Value = 18.57182518074609
$ws.Range("C13").Value = 14.40484932087427
$ws.Range("E13").Value = 21.26412684363698
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 12.88663860233442
$ws.Range("H13").Value = 10.11356929988293
$ws.Range("O13").Value = 13.24264978487028

$ws.Range("B14").Value = 18.423306949444
$ws.Range("C14").Value = 14.30556245204506
$ws.Range("E14").Value = 21.07872347282035
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 12.90211608121778
$ws.Range("H14").Value = 10.12753589933
$ws.Range("O14").Value = 13.26637478474552

$ws.Range("B15").Value = 18.33170842422511
$ws.Range("C15").Value = 14.24438197460351
$ws.Range("E15").Value = 20.96436516790626
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 12.91194365462062
$ws.Range("H15").Value = 10.13614962431112
$ws.Range("O15").Value = 13.28103902409469

$ws.Range("B16").Value = 17.79719947652686
$ws.Range("C16").Value = 13.88823595361151
$ws.Range("E16").Value = 20.29683823289148
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 12.97352493000919
$ws.Range("H16").Value = 10.18641127852554
$ws.Range("O16").Value = 13.3670868793454

$ws.Range("B17").Value = 17.4607498456746
$ws.Range("C17").Value = 13.66485005086138
$ws.Range("E17").Value = 19.87643101304201
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 13.01593559253669
$ws.Range("H17").Value = 10.21804704377929
$ws.Range("O17").Value = 13.42165995264804

$ws.Range("B18").Value = 17.26409970908581
$ws.Range("C18").Value = 13.53458274101079
$ws.Range("E18").Value = 19.63060685712216
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 13.0420056115821
$ws.Range("H18").Value = 10.23653735672922
$ws.Range("O18").Value = 13.45370090799639

$ws.Range("B19").Value = 17.19698001282305
$ws.Range("C19").Value = 13.49017247682757
$ws.Range("E19").Value = 19.54668405087126
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 13.0511179163445
$ws.Range("H19").Value = 10.24284835888306
$ws.Range("O19").Value = 13.46466100904018

$ws.Range("B20").Value = 17.49688994011531
$ws.Range("C20").Value = 13.6888147239608
$ws.Range("E20").Value = 19.92159958113326
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 13.01124679869129
$ws.Range("H20").Value = 10.21464889858692
$ws.Range("O20").Value = 13.41578298605924

$ws.Range("B21").Value = 18.46706101028684
$ws.Range("C21").Value = 14.33480144505942
$ws.Range("E21").Value = 21.13334621781294
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 12.89749748399104
$ws.Range("H21").Value = 10.12342131834522
$ws.Range("O21").Value = 13.25937864590616

$ws.Range("B22").Value = 19.07308669108375
$ws.Range("C22").Value = 14.74073195987815
$ws.Range("E22").Value = 21.88975319768856
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 12.83862519052012
$ws.Range("H22").Value = 10.06642848556156
$ws.Range("O22").Value = 13.16305560810239

$ws.Range("B23").Value = 18.75220989935573
$ws.Range("C23").Value = 14.52558390527435
$ws.Range("E23").Value = 21.48928524794271
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 12.86860574118253
$ws.Range("H23").Value = 10.0966055435417
$ws.Range("O23").Value = 13.21392102462953

$ws.Range("B24").Value = 17.4805610336914
$ws.Range("C24").Value = 13.67798602119129
$ws.Range("E24").Value = 19.90119172772721
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 13.01336135177874
$ws.Range("H24").Value = 10.21618425782632
$ws.Range("O24").Value = 13.41843788845844

$ws.Range("B25").Value = 15.98925896547247
$ws.Range("C25").Value = 12.69594731825377
$ws.Range("E25").Value = 18.03445279330115
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 13.23324388238343
$ws.Range("H25").Value = 10.35641513569769
$ws.Range("O25").Value = 13.66392040037251
